# Automatische test-sync: 2025-07-23 22:31:50
# Adds Testmail #8 ("Planning / Afspraak") as a new row to the Logs sheet,
# rolls the matching count into the Dashboard sheet, and extends the
# Dashboard bar chart's category/value ranges to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 18 with the new test-mail data.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value = "Kun jij een demo inplannen bij Van Dijk deze week?"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Testmail #8: Kun jij een demo inplannen bij Van Dijk deze week?"
$logs.Range("D18").Value = "Planning / Afspraak"
$logs.Range("E18").Value = "Beste afzender,`nBedankt voor je e-mail. Voor het inplannen van een demo bij Van Dijk deze week, vragen we je om het volgende te doen:`n- Stuur ons de beschikbare data en tijden voor de demo.`n- Geef ons wat meer informatie over de gewenste inhoud van de demo, zodat we ons goed kunnen voorbereiden.`nZodra we deze informatie hebben ontvangen, zullen we ons best doen om een geschikte datum en tijd te vinden. `nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F18").Value = "2025-07-23 22:31:45"
$logs.Range("G18").Value = "Ja"
$logs.Range("H18").Value = "Nee"
$logs.Range("I18").Value = "Ja"
$logs.Range("J18").Value = "Nee"

# The multi-line E18 text makes Excel auto-expand the row; AutoFit keeps
# the row height implicit (matching the other data rows, none of which
# carry an explicit ht/customHeight override).
$logs.Rows.Item(18).AutoFit()

# Extend every conditional-formatting rule group on the Logs sheet so it
# keeps covering rows 2-18 (it previously stopped at row 17).
$colLetters = @("D", "G", "H", "I", "J")
foreach ($col in $colLetters) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "17")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "18")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append the matching category/count row.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Planning / Afspraak"
$dash.Range("B8").Value = 1

# ---------------------------------------------------------------------
# 3. Dashboard chart: extend the category + value series ranges to
#    include the new row 8 (leave the series-name reference untouched).
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$8"
$series.Values = "='Dashboard'!`$B`$2:`$B`$8"

Write-Output "Logs row 18 + Dashboard row 8 + chart range added"
